# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Add a new "ODI Bowling Extra" worksheet (after the existing sheets) with
#    MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL scraped data.
# 2) Clean up stray empty cells that were previously written into the
#    "ODI Batting Extra" sheet (columns B-E) where no data was scraped.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "ODI Bowling Extra" sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlExtra.Name = "ODI Bowling Extra"

# Header row
$bowlExtra.Range("A1").Value = "MATCH_CODE"
$bowlExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$header = $bowlExtra.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$rows = @(
    @("4262", "", ""),
    @("4267", "0", "10.00%"),
    @("4391", "0", "20.00%"),
    @("4394", "1", "30.00%"),
    @("4397", "0", "10.00%"),
    @("4426", "", ""),
    @("4427", "0", ""),
    @("4442", "4", "50.00%"),
    @("4444", "2", "10.00%"),
    @("4446", "0", "10.00%"),
    @("4448", "1", "30.00%"),
    @("4466", "3", "10.00%"),
    @("4467", "1", ""),
    @("4468", "2", "30.00%"),
    @("4475", "0", "10.00%"),
    @("4478", "", ""),
    @("4492", "2", "10.00%"),
    @("4496", "0", "20.00%"),
    @("4605", "0", ""),
    @("4608", "0", "20.00%")
)

$r = 2
foreach ($row in $rows) {
    # Leading "'" forces the (numeric-looking) scraped values to be stored
    # as plain text, matching the rest of the workbook's scraped columns.
    $bowlExtra.Cells.Item($r, 1).Value = "'" + $row[0]
    if ($row[1] -ne "") { $bowlExtra.Cells.Item($r, 2).Value = "'" + $row[1] }
    if ($row[2] -ne "") { $bowlExtra.Cells.Item($r, 3).Value = "'" + $row[2] }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Clear stray empty (no scraped value) cells on "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "B5", "C5", "D5", "E5",
    "C12", "D12", "E12",
    "B14", "C14", "D14", "E14",
    "C15", "D15", "E15",
    "B16", "C16", "D16", "E16",
    "B18", "C18", "D18", "E18",
    "C19", "D19", "E19"
)

foreach ($ref in $emptyCells) {
    $battingExtra.Range($ref).ClearContents()
}
